{"js": "// Applies the text-content changes described by the diff to the \"NPCs (non prof.)\" doc.\n// Pure-formatting churn in the diff (w:proofErr gramStart/gramEnd markers that Word's\n// grammar checker stamps around runs, and run merges/splits that leave the visible text\n// identical) does not change the document's text and is not attempted here \u2014 only the\n// edits that change what the document actually says are applied.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"Learning druidcraft...\" \u2014 drop \"which is \" before \"very dangerous work\".\nawait replaceOnce(\n  \"clove farm- which is very dangerous work-\",\n  \"clove farm- very dangerous work-\"\n);\n\n// 2. Diane \u2014 trim the run-on \"attending the academy...\" sentence down to \"She stays...\".\nawait replaceOnce(\n  \"She\\u2019s attending the academy at the same grade level as the party and is staying in one of the dorms next door.\",\n  \"She stays in one of the dorms next door.\"\n);\n\n// 3. Diane \u2014 \"best she could\" -> \"best she can\".\nawait replaceOnce(\"best she could\", \"best she can\");\n\n// 4. Jupo \u2014 expand the hangout-locale list.\nawait replaceOnce(\n  \"hang out with people in popular locales like the plaza, Ondor Ruin, and Hwen\",\n  \"hang out with people in popular locales like the plaza, around the dorms, Hwen, and sometimes Ondor Ruin\"\n);\n\n// 5. Dravik \u2014 fill in the previously-empty bullet with the wrestling-family blurb.\n{\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n  let anchor = -1;\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text.indexOf(\"committed to becoming the best, most powerful spellcaster\") >= 0) {\n      anchor = i;\n      break;\n    }\n  }\n  if (anchor === -1) {\n    throw new Error(\"Could not find Dravik anchor paragraph\");\n  }\n  const target = paras.items[anchor + 1];\n  target.insertText(\n    \"Comes from a family of champion wrestlers in his dwarven clan. Despite him not following their same legacy, they are very supportive of his arcane studies and wish him the best.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 6. Reil \u2014 add \"he's \" before \"often carrying textbooks\".\nawait replaceOnce(\n  \"dragging on the ground, often carrying textbooks\",\n  \"dragging on the ground, he\\u2019s often carrying textbooks\"\n);\n\n// 7. Reil \u2014 remove the standalone \"In one grade level above the party\" bullet.\n{\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text === \"In one grade level above the party\") {\n      paras.items[i].delete();\n      await context.sync();\n      break;\n    }\n  }\n}\n\n// 8. Taffy-cart tiefling \u2014 replace the \"heard tale\" backstory with the \"alumni\" one.\nawait replaceOnce(\n  \". Heard tale of such a place full of eccentric and rich young adults and simply followed the business opportunity. \",\n  \". An alumni saw her capabilities and told her about the campus, but she instead took it up as a business opportunity. \"\n);\n\n// 9. Donna \u2014 \"(non-racial)\" -> \"(non-firbolg)\".\nawait replaceOnce(\"Was never versed in (non-racial) \", \"Was never versed in (non-firbolg) \");\n", "ps1": "# Applies the text-content changes described by the diff to the \"NPCs (non prof.)\" doc.\n# Pure-formatting churn in the diff (w:proofErr gramStart/gramEnd markers that Word's\n# grammar checker stamps around runs, and run merges/splits that leave the visible text\n# identical) does not change the document's text and is not attempted here - only the\n# edits that change what the document actually says are applied.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $searchText, $replaceText) {\n  $find = $doc.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $searchText\n  $find.Replacement.Text = $replaceText\n  $wdFindContinue = 1\n  $wdReplaceOne = 1\n  $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne)\n  if (-not $found) {\n    throw \"Text not found: $searchText\"\n  }\n}\n\n# 1. \"Learning druidcraft...\" - drop \"which is \" before \"very dangerous work\".\nReplace-Text $d \"clove farm- which is very dangerous work-\" \"clove farm- very dangerous work-\"\n\n# 2. Diane - trim the run-on \"attending the academy...\" sentence down to \"She stays...\".\nReplace-Text $d \"She\u2019s attending the academy at the same grade level as the party and is staying in one of the dorms next door.\" \"She stays in one of the dorms next door.\"\n\n# 3. Diane - \"best she could\" -> \"best she can\".\nReplace-Text $d \"best she could\" \"best she can\"\n\n# 4. Jupo - expand the hangout-locale list.\nReplace-Text $d \"hang out with people in popular locales like the plaza, Ondor Ruin, and Hwen\" \"hang out with people in popular locales like the plaza, around the dorms, Hwen, and sometimes Ondor Ruin\"\n\n# 5. Dravik - fill in the previously-empty bullet with the wrestling-family blurb.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.Text -like \"*committed to becoming the best, most powerful spellcaster*\") {\n    $anchorIndex = $i\n    break\n  }\n}\nif ($anchorIndex -eq -1) {\n  throw \"Could not find Dravik anchor paragraph\"\n}\n$target = $d.Paragraphs.Item($anchorIndex + 1)\n$target.Range.InsertBefore(\"Comes from a family of champion wrestlers in his dwarven clan. Despite him not following their same legacy, they are very supportive of his arcane studies and wish him the best.\")\n\n# 6. Reil - add \"he's \" before \"often carrying textbooks\".\nReplace-Text $d \"dragging on the ground, often carrying textbooks\" \"dragging on the ground, he\u2019s often carrying textbooks\"\n\n# 7. Reil - remove the standalone \"In one grade level above the party\" bullet.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text.TrimEnd() -eq \"In one grade level above the party\") {\n    $p.Range.Delete()\n    break\n  }\n}\n\n# 8. Taffy-cart tiefling - replace the \"heard tale\" backstory with the \"alumni\" one.\nReplace-Text $d \". Heard tale of such a place full of eccentric and rich young adults and simply followed the business opportunity. \" \". An alumni saw her capabilities and told her about the campus, but she instead took it up as a business opportunity. \"\n\n# 9. Donna - \"(non-racial)\" -> \"(non-firbolg)\".\nReplace-Text $d \"Was never versed in (non-racial) \" \"Was never versed in (non-firbolg) \"\n"}
